$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rodada 1")

$ws.Range("D2").Value = 48.5
$ws.Range("F2").Value = 43.56

$ws.Range("D3").Value = 46.79
$ws.Range("F3").Value = 66.86

$ws.Range("D4").Value = 65.06
$ws.Range("F4").Value = 47.16

$ws.Range("D5").Value = 43.56
$ws.Range("F5").Value = 59.69

$ws.Range("D7").Value = 71.46
$ws.Range("F7").Value = 57.26

$ws.Range("D8").Value = 71.36
$ws.Range("F8").Value = 49.36

$ws.Range("D9").Value = 50.69
$ws.Range("F9").Value = 63.76

$ws.Range("D10").Value = 56.09
$ws.Range("F10").Value = 81.76

$ws.Range("D11").Value = 50.76
$ws.Range("F11").Value = 59.65

$ws.Range("D12").Value = 33.96
$ws.Range("F12").Value = 58.26

$ws.Range("D13").Value = 52.39
$ws.Range("F13").Value = 54.95

$ws.Range("D14").Value = 48.89
$ws.Range("F14").Value = 59.36

$ws.Range("D15").Value = 71.16
$ws.Range("F15").Value = 48.29

$ws.Range("D16").Value = 38.66

$ws.Range("D17").Value = 63.76
$ws.Range("F17").Value = 54.66
